$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Output ($win.SheetViews | Get-Member | Select-String "Name=")
